$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 18 with the new entry for 19/5/2025
$ws.Range("D18").Value = "19/5/2025"
$ws.Range("E18").Value = 135
$ws.Range("F18").Value = 218
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 650
$ws.Range("J18").Value = "N/A"

# Update the active cell selection to match the author's final cursor position
$ws.Range("J21").Select()
